$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues($ws, $RowNum, $Values) {
    $n = $Values.Length
    $arr = New-Object 'object[,]' 1, $n
    for ($i = 0; $i -lt $n; $i++) {
        $arr[0, $i] = $Values[$i]
    }
    $startCell = $ws.Cells.Item($RowNum, 1)
    $endCell = $ws.Cells.Item($RowNum, $n)
    $rng = $ws.Range($startCell, $endCell)
    $rng.Value = $arr
}

# --- Un-merge the old grouped headers (Tackles / Challenges / Blocks spans) ---
# Row 1 now carries a real column header in every cell instead of a merged group label.
$ws.Range("H1:L1").UnMerge()
$ws.Range("M1:P1").UnMerge()
$ws.Range("Q1:S1").UnMerge()

# --- Row 1: real per-column header labels ---
$row1 = @('Player ID', 'Player', '#', 'Nation', 'Pos', 'Age', '90s', 'Tkl', 'TklW', 'Def 3rd', 'Mid 3rd', 'Att 3rd', 'Cha', 'Att', 'Tkl%', 'Lost', 'Blocks', 'Sh', 'Pass', 'Int', 'Tkl+Int', 'Clr', 'Err')
Set-RowValues $ws 1 $row1

# --- Row 2: secondary header row (kept, but now hidden) ---
$row2 = @($null, 'Player', '#', 'Nation', 'Pos', 'Age', 'Min', 'Tkl', 'TklW', 'Def 3rd', 'Mid 3rd', 'Att 3rd', 'Tkl', 'Att', 'Tkl%', 'Lost', 'Blocks', 'Sh', 'Pass', 'Int', 'Tkl+Int', 'Clr', 'Err')
Set-RowValues $ws 2 $row2

# --- Row 3: stays blank, will be hidden below ---

# --- Player data rows ---
Set-RowValues $ws 4 @(0, 'Klauss', 9, 'br BRA', 'FW', '26-010', 90, 2, 0, 1, 1, 0, 2, 3, 66.7, 1, 1, 0, 1, 0, 2, 0, 0)
Set-RowValues $ws 5 @(1, 'Jared Stroud', 8, 'us USA', 'LW', '26-244', 65, 0, 0, 0, 0, 0, 0, 1, 0, 1, 0, 0, 0, 1, 1, 0, 0)
Set-RowValues $ws 6 @(2, 'Célio Pompeu', 12, 'br BRA', 'LW,LB', '23-091', 25, 1, 1, 0, 1, 0, 0, 4, 0, 4, 2, 0, 2, 0, 1, 0, 0)
Set-RowValues $ws 7 @(3, 'Rasmus Alm', 21, 'se SWE', 'RW', '27-206', 71, 3, 3, 0, 2, 1, 0, 0, 0, 0, 2, 0, 2, 0, 3, 0, 0)
Set-RowValues $ws 8 @(4, 'Nicholas Gioacchini', 11, 'us USA', 'RW,FW', '22-229', 19, 1, 0, 0, 1, 0, 0, 0, 0, 0, 0, 0, 0, 0, 1, 0, 0)
Set-RowValues $ws 9 @(5, 'Eduard Löwen', 10, 'de GER', 'AM,CM', '26-042', 90, 1, 0, 1, 0, 0, 1, 1, 100, 0, 3, 0, 3, 1, 2, 1, 0)
Set-RowValues $ws 10 @(6, 'Miguel Perez', 28, 'us USA', 'DM', '17-317', 65, 3, 2, 3, 0, 0, 3, 5, 60, 2, 1, 0, 1, 2, 5, 1, 0)
Set-RowValues $ws 11 @(7, 'Tomáš Ostrák', 7, 'cz CZE', 'DM,CM', '23-034', 25, 1, 1, 1, 0, 0, 1, 3, 33.3, 2, 0, 0, 0, 0, 1, 0, 0)
Set-RowValues $ws 12 @(8, 'Indiana Vassilev', 19, 'us USA', 'DM', '22-023', 71, 4, 3, 1, 3, 0, 2, 3, 66.7, 1, 1, 0, 1, 0, 4, 0, 0)
Set-RowValues $ws 13 @(9, 'Akil Watts', 20, 'us USA', 'DM,CM', '23-035', 19, 2, 1, 0, 2, 0, 1, 1, 100, 0, 0, 0, 0, 1, 3, 0, 0)
Set-RowValues $ws 14 @(10, 'John Nelson', 14, 'us USA', 'LB', '24-243', 85, 5, 4, 2, 3, 0, 3, 4, 75, 1, 1, 0, 1, 2, 7, 6, 0)
Set-RowValues $ws 15 @(11, 'Lucas Bartlett', 24, 'us USA', 'CB', '25-228', 5, 0, 0, 0, 0, 0, 0, 0, 0, 0, 1, 0, 1, 2, 2, 1, 0)
Set-RowValues $ws 16 @(12, 'Kyle Hiebert', 22, 'ca CAN', 'CB', '25-224', 90, 0, 0, 0, 0, 0, 0, 3, 0, 3, 2, 1, 1, 3, 3, 5, 0)
Set-RowValues $ws 17 @(13, 'Tim Parker', 26, 'us USA', 'CB', '30-016', 90, 1, 0, 1, 0, 0, 1, 1, 100, 0, 0, 0, 0, 1, 2, 8, 0)
Set-RowValues $ws 18 @(14, 'Jake Nerwinski', 2, 'us USA', 'RB', '28-145', 90, 2, 2, 1, 0, 1, 1, 2, 50, 1, 3, 0, 3, 4, 6, 1, 0)
Set-RowValues $ws 19 @(15, 'Roman Bürki', 1, 'ch SUI', 'GK', '32-117', 90, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
Set-RowValues $ws 20 @(16, '16 Players', $null, $null, $null, $null, 990, 26, 17, 11, 13, 2, 15, 31, 48.4, 16, 17, 1, 16, 17, 43, 23, 0)

# --- Hide the helper/summary rows (sub-header, spacer, and the "16 Players" totals row) ---
$ws.Rows.Item(2).Hidden = $true
$ws.Rows.Item(3).Hidden = $true
$ws.Rows.Item(20).Hidden = $true

# --- Restore the active selection used when the file was last saved ---
$ws.Range("O21").Select()

Write-Host "Defensive actions sheet cleaned up."
